$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing recalculated statistic values (rows 2-163) ---
$ws.Cells.Item(2,7).Value = 1.98507356735162
$ws.Cells.Item(2,8).Value = 6.16845766552593
$ws.Cells.Item(2,9).Value = 3.7915
$ws.Cells.Item(9,7).Value = 0.0400494314567304
$ws.Cells.Item(10,7).Value = 0.0400494314567304
$ws.Cells.Item(19,7).Value = 1.84382203469846
$ws.Cells.Item(19,8).Value = 6.16845766552593
$ws.Cells.Item(19,9).Value = 3.481
$ws.Cells.Item(26,7).Value = 0.0376223336911737
$ws.Cells.Item(27,7).Value = 0.0376223336911737
$ws.Cells.Item(43,7).Value = 0.0298831484335109
$ws.Cells.Item(44,7).Value = 0.0298831484335109
$ws.Cells.Item(45,7).Value = 0.569659766834931
$ws.Cells.Item(46,7).Value = 0.569659766834931
$ws.Cells.Item(47,7).Value = 0.612811521225452
$ws.Cells.Item(48,7).Value = 0.612811521225452
$ws.Cells.Item(60,7).Value = 0.0267827050047322
$ws.Cells.Item(61,7).Value = 0.0267827050047322
$ws.Cells.Item(62,7).Value = 0.576865450696221
$ws.Cells.Item(63,7).Value = 0.576865450696221
$ws.Cells.Item(64,7).Value = 0.616768034734456
$ws.Cells.Item(65,7).Value = 0.616768034734456
$ws.Cells.Item(77,7).Value = 0.0182841968096218
$ws.Cells.Item(77,12).Value = 0.00476
$ws.Cells.Item(78,7).Value = 0.0182841968096218
$ws.Cells.Item(78,12).Value = 0.00476
$ws.Cells.Item(79,7).Value = 0.6314987840295549
$ws.Cells.Item(80,7).Value = 0.6314987840295549
$ws.Cells.Item(81,7).Value = 0.670001368067789
$ws.Cells.Item(82,7).Value = 0.670001368067789
$ws.Cells.Item(94,7).Value = 0.0178490434255875
$ws.Cells.Item(94,12).Value = 0.00476
$ws.Cells.Item(95,7).Value = 0.0178490434255875
$ws.Cells.Item(95,12).Value = 0.00476
$ws.Cells.Item(96,7).Value = 0.683503784029555
$ws.Cells.Item(97,7).Value = 0.683503784029555
$ws.Cells.Item(98,7).Value = 0.720834701401122
$ws.Cells.Item(99,7).Value = 0.720834701401122
$ws.Cells.Item(104,7).Value = 1.25621221461009
$ws.Cells.Item(111,7).Value = 0.0195000112240748
$ws.Cells.Item(111,12).Value = 0.0074
$ws.Cells.Item(112,7).Value = 0.0195000112240748
$ws.Cells.Item(112,12).Value = 0.0074
$ws.Cells.Item(113,7).Value = 0.679602117362888
$ws.Cells.Item(114,7).Value = 0.679602117362888
$ws.Cells.Item(115,7).Value = 0.720318034734456
$ws.Cells.Item(116,7).Value = 0.720318034734456
$ws.Cells.Item(121,7).Value = 1.21170548579213
$ws.Cells.Item(128,7).Value = 0.0214605638528465
$ws.Cells.Item(129,7).Value = 0.0214605638528465
$ws.Cells.Item(130,7).Value = 0.689152555906988
$ws.Cells.Item(131,7).Value = 0.689152555906988
$ws.Cells.Item(132,7).Value = 0.732054276809587
$ws.Cells.Item(133,7).Value = 0.732054276809587
$ws.Cells.Item(138,7).Value = 1.02249140920798
$ws.Cells.Item(145,7).Value = 0.0210105708640519
$ws.Cells.Item(146,7).Value = 0.0210105708640519
$ws.Cells.Item(155,7).Value = 0.923802405819382
$ws.Cells.Item(158,7).Value = 1947.68795611256
$ws.Cells.Item(158,8).Value = 12247.5894106411
$ws.Cells.Item(159,7).Value = 1947.68795611256
$ws.Cells.Item(159,8).Value = 12247.5894106411
$ws.Cells.Item(160,7).Value = 1947.68795611256
$ws.Cells.Item(160,8).Value = 12247.5894106411
$ws.Cells.Item(161,7).Value = 1947.68795611256
$ws.Cells.Item(161,8).Value = 12247.5894106411
$ws.Cells.Item(162,7).Value = 0.0194721850567072
$ws.Cells.Item(163,7).Value = 0.0194721850567072

# --- Append new rows 172-188 for 2019-2023 year range ---
$ws.Cells.Item(172,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(172,2).Value = "Visual Clarity (Sediment class 4)"
$ws.Cells.Item(172,3).Value = "D"
$ws.Cells.Item(172,4).Value = "2019 - 2023"
$ws.Cells.Item(172,5).Value = "RepSite"
$ws.Cells.Item(172,6).Value = 0.64
$ws.Cells.Item(172,7).Value = 0.839355771552214
$ws.Cells.Item(172,8).Value = 2.538
$ws.Cells.Item(172,9).Value = 2.09
$ws.Cells.Item(172,12).Value = 1.29
$ws.Cells.Item(172,13).Value = 1.744
$ws.Cells.Item(172,14).Value = 2.0192
$ws.Cells.Item(172,15).Value = 1863911.8
$ws.Cells.Item(172,16).Value = 5542684.5
$ws.Cells.Item(172,17).Value = "Tararua District"
$ws.Cells.Item(172,18).Value = "Manawatū"
$ws.Cells.Item(172,19).Value = "Weber - Tamaki"
$ws.Cells.Item(172,20).Value = "Mana_2b"
$ws.Cells.Item(172,21).Value = "m"
$ws.Cells.Item(173,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(173,2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(173,3).Value = "D"
$ws.Cells.Item(173,4).Value = "2019 - 2023"
$ws.Cells.Item(173,5).Value = "RepSite"
$ws.Cells.Item(173,6).Value = 0.047
$ws.Cells.Item(173,7).Value = 0.060728813559322
$ws.Cells.Item(173,8).Value = 0.274
$ws.Cells.Item(173,9).Value = 0.1125
$ws.Cells.Item(173,12).Value = 0.0585
$ws.Cells.Item(173,13).Value = 0.08846999999999999
$ws.Cells.Item(173,14).Value = 0.10256
$ws.Cells.Item(173,15).Value = 1863911.8
$ws.Cells.Item(173,16).Value = 5542684.5
$ws.Cells.Item(173,17).Value = "Tararua District"
$ws.Cells.Item(173,18).Value = "Manawatū"
$ws.Cells.Item(173,19).Value = "Weber - Tamaki"
$ws.Cells.Item(173,20).Value = "Mana_2b"
$ws.Cells.Item(173,21).Value = "mg/L"
$ws.Cells.Item(174,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(174,2).Value = "DRP (Median)"
$ws.Cells.Item(174,3).Value = "D"
$ws.Cells.Item(174,4).Value = "2019 - 2023"
$ws.Cells.Item(174,5).Value = "RepSite"
$ws.Cells.Item(174,6).Value = 0.047
$ws.Cells.Item(174,7).Value = 0.060728813559322
$ws.Cells.Item(174,8).Value = 0.274
$ws.Cells.Item(174,9).Value = 0.1125
$ws.Cells.Item(174,12).Value = 0.0585
$ws.Cells.Item(174,13).Value = 0.08846999999999999
$ws.Cells.Item(174,14).Value = 0.10256
$ws.Cells.Item(174,15).Value = 1863911.8
$ws.Cells.Item(174,16).Value = 5542684.5
$ws.Cells.Item(174,17).Value = "Tararua District"
$ws.Cells.Item(174,18).Value = "Manawatū"
$ws.Cells.Item(174,19).Value = "Weber - Tamaki"
$ws.Cells.Item(174,20).Value = "Mana_2b"
$ws.Cells.Item(174,21).Value = "mg/L"
$ws.Cells.Item(175,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(175,2).Value = "E coli (>260)"
$ws.Cells.Item(175,3).Value = "E"
$ws.Cells.Item(175,4).Value = "2019 - 2023"
$ws.Cells.Item(175,5).Value = "RepSite"
$ws.Cells.Item(175,6).Value = 790
$ws.Cells.Item(175,7).Value = 1744.63710865493
$ws.Cells.Item(175,8).Value = 12247.5894106411
$ws.Cells.Item(175,9).Value = 8877.9
$ws.Cells.Item(175,10).Value = 61.0169491525424
$ws.Cells.Item(175,11).Value = 89.83050847457631
$ws.Cells.Item(175,12).Value = 1050
$ws.Cells.Item(175,13).Value = 2100
$ws.Cells.Item(175,14).Value = 6625.88
$ws.Cells.Item(175,15).Value = 1863911.8
$ws.Cells.Item(175,16).Value = 5542684.5
$ws.Cells.Item(175,17).Value = "Tararua District"
$ws.Cells.Item(175,18).Value = "Manawatū"
$ws.Cells.Item(175,19).Value = "Weber - Tamaki"
$ws.Cells.Item(175,20).Value = "Mana_2b"
$ws.Cells.Item(175,21).Value = "% exceedances over 260/100 mL"
$ws.Cells.Item(176,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(176,2).Value = "E coli (>540)"
$ws.Cells.Item(176,3).Value = "E"
$ws.Cells.Item(176,4).Value = "2019 - 2023"
$ws.Cells.Item(176,5).Value = "RepSite"
$ws.Cells.Item(176,6).Value = 790
$ws.Cells.Item(176,7).Value = 1744.63710865493
$ws.Cells.Item(176,8).Value = 12247.5894106411
$ws.Cells.Item(176,9).Value = 8877.9
$ws.Cells.Item(176,10).Value = 61.0169491525424
$ws.Cells.Item(176,11).Value = 89.83050847457631
$ws.Cells.Item(176,12).Value = 1050
$ws.Cells.Item(176,13).Value = 2100
$ws.Cells.Item(176,14).Value = 6625.88
$ws.Cells.Item(176,15).Value = 1863911.8
$ws.Cells.Item(176,16).Value = 5542684.5
$ws.Cells.Item(176,17).Value = "Tararua District"
$ws.Cells.Item(176,18).Value = "Manawatū"
$ws.Cells.Item(176,19).Value = "Weber - Tamaki"
$ws.Cells.Item(176,20).Value = "Mana_2b"
$ws.Cells.Item(176,21).Value = "% exceedances over 540/100 mL"
$ws.Cells.Item(177,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(177,2).Value = "E coli (Median)"
$ws.Cells.Item(177,3).Value = "E"
$ws.Cells.Item(177,4).Value = "2019 - 2023"
$ws.Cells.Item(177,5).Value = "RepSite"
$ws.Cells.Item(177,6).Value = 790
$ws.Cells.Item(177,7).Value = 1744.63710865493
$ws.Cells.Item(177,8).Value = 12247.5894106411
$ws.Cells.Item(177,9).Value = 8877.9
$ws.Cells.Item(177,10).Value = 61.0169491525424
$ws.Cells.Item(177,11).Value = 89.83050847457631
$ws.Cells.Item(177,12).Value = 1050
$ws.Cells.Item(177,13).Value = 2100
$ws.Cells.Item(177,14).Value = 6625.88
$ws.Cells.Item(177,15).Value = 1863911.8
$ws.Cells.Item(177,16).Value = 5542684.5
$ws.Cells.Item(177,17).Value = "Tararua District"
$ws.Cells.Item(177,18).Value = "Manawatū"
$ws.Cells.Item(177,19).Value = "Weber - Tamaki"
$ws.Cells.Item(177,20).Value = "Mana_2b"
$ws.Cells.Item(177,21).Value = "E. coli/100 mL"
$ws.Cells.Item(178,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(178,2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(178,3).Value = "E"
$ws.Cells.Item(178,4).Value = "2019 - 2023"
$ws.Cells.Item(178,5).Value = "RepSite"
$ws.Cells.Item(178,6).Value = 790
$ws.Cells.Item(178,7).Value = 1744.63710865493
$ws.Cells.Item(178,8).Value = 12247.5894106411
$ws.Cells.Item(178,9).Value = 8877.9
$ws.Cells.Item(178,10).Value = 61.0169491525424
$ws.Cells.Item(178,11).Value = 89.83050847457631
$ws.Cells.Item(178,12).Value = 1050
$ws.Cells.Item(178,13).Value = 2100
$ws.Cells.Item(178,14).Value = 6625.88
$ws.Cells.Item(178,15).Value = 1863911.8
$ws.Cells.Item(178,16).Value = 5542684.5
$ws.Cells.Item(178,17).Value = "Tararua District"
$ws.Cells.Item(178,18).Value = "Manawatū"
$ws.Cells.Item(178,19).Value = "Weber - Tamaki"
$ws.Cells.Item(178,20).Value = "Mana_2b"
$ws.Cells.Item(178,21).Value = "E. coli/100 mL"
$ws.Cells.Item(179,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(179,2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(179,3).Value = "A"
$ws.Cells.Item(179,4).Value = "2019 - 2023"
$ws.Cells.Item(179,5).Value = "RepSite"
$ws.Cells.Item(179,6).Value = 0.01361
$ws.Cells.Item(179,7).Value = 0.0170344305837054
$ws.Cells.Item(179,8).Value = 0.169853842472124
$ws.Cells.Item(179,9).Value = 0.03454
$ws.Cells.Item(179,12).Value = 0.01169
$ws.Cells.Item(179,13).Value = 0.02363
$ws.Cells.Item(179,14).Value = 0.03208
$ws.Cells.Item(179,15).Value = 1863911.8
$ws.Cells.Item(179,16).Value = 5542684.5
$ws.Cells.Item(179,17).Value = "Tararua District"
$ws.Cells.Item(179,18).Value = "Manawatū"
$ws.Cells.Item(179,19).Value = "Weber - Tamaki"
$ws.Cells.Item(179,20).Value = "Mana_2b"
$ws.Cells.Item(179,21).Value = "mg NH4-N/L"
$ws.Cells.Item(180,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(180,2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(180,3).Value = "A"
$ws.Cells.Item(180,4).Value = "2019 - 2023"
$ws.Cells.Item(180,5).Value = "RepSite"
$ws.Cells.Item(180,6).Value = 0.01361
$ws.Cells.Item(180,7).Value = 0.0170344305837054
$ws.Cells.Item(180,8).Value = 0.169853842472124
$ws.Cells.Item(180,9).Value = 0.03454
$ws.Cells.Item(180,12).Value = 0.01169
$ws.Cells.Item(180,13).Value = 0.02363
$ws.Cells.Item(180,14).Value = 0.03208
$ws.Cells.Item(180,15).Value = 1863911.8
$ws.Cells.Item(180,16).Value = 5542684.5
$ws.Cells.Item(180,17).Value = "Tararua District"
$ws.Cells.Item(180,18).Value = "Manawatū"
$ws.Cells.Item(180,19).Value = "Weber - Tamaki"
$ws.Cells.Item(180,20).Value = "Mana_2b"
$ws.Cells.Item(180,21).Value = "mg NH4-N/L"
$ws.Cells.Item(181,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(181,2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(181,3).Value = "A"
$ws.Cells.Item(181,4).Value = "2019 - 2023"
$ws.Cells.Item(181,5).Value = "RepSite"
$ws.Cells.Item(181,6).Value = 0.589
$ws.Cells.Item(181,7).Value = 0.690779661016949
$ws.Cells.Item(181,8).Value = 1.68
$ws.Cells.Item(181,9).Value = 1.4885
$ws.Cells.Item(181,12).Value = 0.259
$ws.Cells.Item(181,13).Value = 1.2447
$ws.Cells.Item(181,14).Value = 1.3634
$ws.Cells.Item(181,15).Value = 1863911.8
$ws.Cells.Item(181,16).Value = 5542684.5
$ws.Cells.Item(181,17).Value = "Tararua District"
$ws.Cells.Item(181,18).Value = "Manawatū"
$ws.Cells.Item(181,19).Value = "Weber - Tamaki"
$ws.Cells.Item(181,20).Value = "Mana_2b"
$ws.Cells.Item(181,21).Value = "mg NO3-N/L"
$ws.Cells.Item(182,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(182,2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(182,3).Value = "A"
$ws.Cells.Item(182,4).Value = "2019 - 2023"
$ws.Cells.Item(182,5).Value = "RepSite"
$ws.Cells.Item(182,6).Value = 0.589
$ws.Cells.Item(182,7).Value = 0.690779661016949
$ws.Cells.Item(182,8).Value = 1.68
$ws.Cells.Item(182,9).Value = 1.4885
$ws.Cells.Item(182,12).Value = 0.259
$ws.Cells.Item(182,13).Value = 1.2447
$ws.Cells.Item(182,14).Value = 1.3634
$ws.Cells.Item(182,15).Value = 1863911.8
$ws.Cells.Item(182,16).Value = 5542684.5
$ws.Cells.Item(182,17).Value = "Tararua District"
$ws.Cells.Item(182,18).Value = "Manawatū"
$ws.Cells.Item(182,19).Value = "Weber - Tamaki"
$ws.Cells.Item(182,20).Value = "Mana_2b"
$ws.Cells.Item(182,21).Value = "mg NO3-N/L"
$ws.Cells.Item(183,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(183,2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(183,4).Value = "2019 - 2023"
$ws.Cells.Item(183,5).Value = "RepSite"
$ws.Cells.Item(183,6).Value = 0.628
$ws.Cells.Item(183,7).Value = 0.727406779661017
$ws.Cells.Item(183,8).Value = 1.759
$ws.Cells.Item(183,9).Value = 1.57715
$ws.Cells.Item(183,12).Value = 0.2885
$ws.Cells.Item(183,13).Value = 1.30311
$ws.Cells.Item(183,14).Value = 1.48224
$ws.Cells.Item(183,15).Value = 1863911.8
$ws.Cells.Item(183,16).Value = 5542684.5
$ws.Cells.Item(183,17).Value = "Tararua District"
$ws.Cells.Item(183,18).Value = "Manawatū"
$ws.Cells.Item(183,19).Value = "Weber - Tamaki"
$ws.Cells.Item(183,20).Value = "Mana_2b"
$ws.Cells.Item(183,21).Value = "g/m3"
$ws.Cells.Item(184,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(184,2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(184,4).Value = "2019 - 2023"
$ws.Cells.Item(184,5).Value = "RepSite"
$ws.Cells.Item(184,6).Value = 0.628
$ws.Cells.Item(184,7).Value = 0.727406779661017
$ws.Cells.Item(184,8).Value = 1.759
$ws.Cells.Item(184,9).Value = 1.57715
$ws.Cells.Item(184,12).Value = 0.2885
$ws.Cells.Item(184,13).Value = 1.30311
$ws.Cells.Item(184,14).Value = 1.48224
$ws.Cells.Item(184,15).Value = 1863911.8
$ws.Cells.Item(184,16).Value = 5542684.5
$ws.Cells.Item(184,17).Value = "Tararua District"
$ws.Cells.Item(184,18).Value = "Manawatū"
$ws.Cells.Item(184,19).Value = "Weber - Tamaki"
$ws.Cells.Item(184,20).Value = "Mana_2b"
$ws.Cells.Item(184,21).Value = "g/m3"
$ws.Cells.Item(185,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(185,2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(185,4).Value = "2019 - 2023"
$ws.Cells.Item(185,5).Value = "RepSite"
$ws.Cells.Item(185,6).Value = 0.91
$ws.Cells.Item(185,7).Value = 1.02440677966102
$ws.Cells.Item(185,8).Value = 3.11
$ws.Cells.Item(185,9).Value = 2.046
$ws.Cells.Item(185,12).Value = 0.53
$ws.Cells.Item(185,13).Value = 1.5835
$ws.Cells.Item(185,14).Value = 1.839
$ws.Cells.Item(185,15).Value = 1863911.8
$ws.Cells.Item(185,16).Value = 5542684.5
$ws.Cells.Item(185,17).Value = "Tararua District"
$ws.Cells.Item(185,18).Value = "Manawatū"
$ws.Cells.Item(185,19).Value = "Weber - Tamaki"
$ws.Cells.Item(185,20).Value = "Mana_2b"
$ws.Cells.Item(185,21).Value = "g/m3"
$ws.Cells.Item(186,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(186,2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(186,4).Value = "2019 - 2023"
$ws.Cells.Item(186,5).Value = "RepSite"
$ws.Cells.Item(186,6).Value = 0.91
$ws.Cells.Item(186,7).Value = 1.02440677966102
$ws.Cells.Item(186,8).Value = 3.11
$ws.Cells.Item(186,9).Value = 2.046
$ws.Cells.Item(186,12).Value = 0.53
$ws.Cells.Item(186,13).Value = 1.5835
$ws.Cells.Item(186,14).Value = 1.839
$ws.Cells.Item(186,15).Value = 1863911.8
$ws.Cells.Item(186,16).Value = 5542684.5
$ws.Cells.Item(186,17).Value = "Tararua District"
$ws.Cells.Item(186,18).Value = "Manawatū"
$ws.Cells.Item(186,19).Value = "Weber - Tamaki"
$ws.Cells.Item(186,20).Value = "Mana_2b"
$ws.Cells.Item(186,21).Value = "g/m3"
$ws.Cells.Item(187,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(187,2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(187,4).Value = "2019 - 2023"
$ws.Cells.Item(187,5).Value = "RepSite"
$ws.Cells.Item(187,6).Value = 0.089
$ws.Cells.Item(187,7).Value = 0.126898305084746
$ws.Cells.Item(187,8).Value = 1.18
$ws.Cells.Item(187,9).Value = 0.345
$ws.Cells.Item(187,12).Value = 0.1
$ws.Cells.Item(187,13).Value = 0.13329
$ws.Cells.Item(187,14).Value = 0.17236
$ws.Cells.Item(187,15).Value = 1863911.8
$ws.Cells.Item(187,16).Value = 5542684.5
$ws.Cells.Item(187,17).Value = "Tararua District"
$ws.Cells.Item(187,18).Value = "Manawatū"
$ws.Cells.Item(187,19).Value = "Weber - Tamaki"
$ws.Cells.Item(187,20).Value = "Mana_2b"
$ws.Cells.Item(187,21).Value = "g/m3"
$ws.Cells.Item(188,1).Value = "Mangatera at u/s T.D.C. Ox Ponds"
$ws.Cells.Item(188,2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(188,4).Value = "2019 - 2023"
$ws.Cells.Item(188,5).Value = "RepSite"
$ws.Cells.Item(188,6).Value = 0.089
$ws.Cells.Item(188,7).Value = 0.126898305084746
$ws.Cells.Item(188,8).Value = 1.18
$ws.Cells.Item(188,9).Value = 0.345
$ws.Cells.Item(188,12).Value = 0.1
$ws.Cells.Item(188,13).Value = 0.13329
$ws.Cells.Item(188,14).Value = 0.17236
$ws.Cells.Item(188,15).Value = 1863911.8
$ws.Cells.Item(188,16).Value = 5542684.5
$ws.Cells.Item(188,17).Value = "Tararua District"
$ws.Cells.Item(188,18).Value = "Manawatū"
$ws.Cells.Item(188,19).Value = "Weber - Tamaki"
$ws.Cells.Item(188,20).Value = "Mana_2b"
$ws.Cells.Item(188,21).Value = "g/m3"